$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename header row (row 1) -------------------------------------------
# Columns A:J were "<name>_old" -> "<name>_FV2410"
# Column K ("diff") is unchanged
# Columns L:U were "<name>_new" -> "<name>_FV2504"
$baseNames = @("Segmentname", "Segmentgruppe", "Segment", "Datenelement", "Segment ID", "Code", "Qualifier", "Beschreibung", "Bedingungsausdruck", "Bedingung")

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 1
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2410"
}

$ws.Cells.Item(1, 11).Value = "diff"

for ($i = 0; $i -lt $baseNames.Count; $i++) {
    $col = $i + 12
    $ws.Cells.Item(1, $col).Value = $baseNames[$i] + "_FV2504"
}

# --- Turn the data range into a table (ListObject) -----------------------
$lo = $ws.ListObjects.Add([Microsoft.Office.Interop.Excel.XlListObjectSourceType]::xlSrcRange, $ws.Range("A1:U80"), $null, [Microsoft.Office.Interop.Excel.XlYesNoGuess]::xlYes)
$lo.Name = "Table1"
$lo.TableStyle = ""
$lo.ShowTableStyleRowStripes = $true
$lo.ShowTableStyleColumnStripes = $false
$lo.ShowTableStyleFirstColumn = $false
$lo.ShowTableStyleLastColumn = $false

# --- Freeze the header row ------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
